$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header date text (stored as plain text string, not an actual date)
$ws.Range("B1").Value = "21/03/2023"

# Row 2 - AMM
$ws.Range("B2").Value = 707.2
$ws.Range("C2").Value = 10

# Row 3 - IPR
$ws.Range("B3").Value = 81
$ws.Range("C3").Value = 10

# Row 4 - MIG
$ws.Range("B4").Value = 120
$ws.Range("C4").Value = 10

# Row 5 - MOB
$ws.Range("B5").Value = 504
$ws.Range("C5").Value = 10

# Row 6 - MOB PRE
$ws.Range("B6").Value = 364
$ws.Range("C6").Value = 10

# Row 7 - MSK
$ws.Range("B7").Value = 110
$ws.Range("C7").Value = 10

# Row 8 - NOT
$ws.Range("B8").Value = 119
$ws.Range("C8").Value = 10

# Row 9 - TEC
$ws.Range("B9").Value = 441
$ws.Range("C9").Value = 10

# Row 10 - TST
$ws.Range("B10").Value = 50
$ws.Range("C10").Value = 10

# Row 11 - VIP (B11 unchanged)
$ws.Range("C11").Value = 10

# Row 12 - WLC
$ws.Range("B12").Value = 44
$ws.Range("C12").Value = 10
